$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(55, 8).Value = 355.9091  # ALC!H55: 384.5 -> 355.9091
$ws.Cells.Item(55, 10).Value = 73  # ALC!J55: 73.75 -> 73
$ws.Cells.Item(55, 12).Value = 73  # ALC!L55: 73.75 -> 73
$ws.Cells.Item(55, 14).Value = -501  # ALC!N55: -501.75 -> -501

$ws.Cells.Item(70, 8).Value = 3181.8333  # ALC!H70: 41669396 -> 3181.8333
$ws.Cells.Item(70, 9).Value = 3566.6667  # ALC!I70: 62502676 -> 3566.6667
$ws.Cells.Item(70, 10).Value = 3053.5557  # ALC!J70: 33336082 -> 3053.5557
$ws.Cells.Item(70, 11).Value = 10700.0001  # ALC!K70: 187508028 -> 10700.0001
$ws.Cells.Item(70, 12).Value = 9160.667099999999  # ALC!L70: 100008246 -> 9160.667099999999
$ws.Cells.Item(70, 13).Value = -10430.0001  # ALC!M70: -187507758 -> -10430.0001
$ws.Cells.Item(70, 14).Value = -9700.667099999999  # ALC!N70: -100008786 -> -9700.667099999999

$ws.Cells.Item(73, 8).Value = 3181.8333  # ALC!H73: 41669396 -> 3181.8333
$ws.Cells.Item(73, 9).Value = 3566.6667  # ALC!I73: 62502676 -> 3566.6667
$ws.Cells.Item(73, 10).Value = 3053.5557  # ALC!J73: 33336082 -> 3053.5557
$ws.Cells.Item(73, 11).Value = 10700.0001  # ALC!K73: 187508028 -> 10700.0001
$ws.Cells.Item(73, 12).Value = 9160.667099999999  # ALC!L73: 100008246 -> 9160.667099999999
$ws.Cells.Item(73, 13).Value = -9764.000100000001  # ALC!M73: -187507092 -> -9764.000100000001
$ws.Cells.Item(73, 14).Value = -11032.6671  # ALC!N73: -100010118 -> -11032.6671

$ws.Cells.Item(75, 8).Value = 0  # ALC!H75: 42500 -> 0
$ws.Cells.Item(75, 10).Value = 0  # ALC!J75: 42500 -> 0
$ws.Cells.Item(75, 12).Value = 0  # ALC!L75: 42500 -> 0
$ws.Cells.Item(75, 14).ClearContents()  # ALC!N75: -44372 -> (removed)

$ws.Cells.Item(78, 8).Value = 0  # ALC!H78: 42500 -> 0
$ws.Cells.Item(78, 10).Value = 0  # ALC!J78: 42500 -> 0
$ws.Cells.Item(78, 12).Value = 0  # ALC!L78: 127500 -> 0
$ws.Cells.Item(78, 14).ClearContents()  # ALC!N78: -136860 -> (removed)

$ws.Cells.Item(80, 8).Value = 30743.883  # ALC!H80: 34776.465 -> 30743.883
$ws.Cells.Item(80, 9).Value = 13096.75  # ALC!I80: 14920.143 -> 13096.75
$ws.Cells.Item(80, 10).Value = 46430.223  # ALC!J80: 52150.75 -> 46430.223
$ws.Cells.Item(80, 11).Value = 39290.25  # ALC!K80: 44760.429 -> 39290.25
$ws.Cells.Item(80, 12).Value = 139290.669  # ALC!L80: 156452.25 -> 139290.669
$ws.Cells.Item(80, 13).Value = -38292.25  # ALC!M80: -43762.429 -> -38292.25
$ws.Cells.Item(80, 14).Value = -141286.669  # ALC!N80: -158448.25 -> -141286.669

$ws.Cells.Item(82, 8).Value = 799.6667  # ALC!H82: 2000 -> 799.6667
$ws.Cells.Item(82, 9).Value = 799.6667  # ALC!I82: 2000 -> 799.6667
$ws.Cells.Item(82, 11).Value = 2399.0001  # ALC!K82: 6000 -> 2399.0001
$ws.Cells.Item(82, 13).Value = -1993.0001  # ALC!M82: -5594 -> -1993.0001

$ws.Cells.Item(83, 8).Value = 30743.883  # ALC!H83: 34776.465 -> 30743.883
$ws.Cells.Item(83, 9).Value = 13096.75  # ALC!I83: 14920.143 -> 13096.75
$ws.Cells.Item(83, 10).Value = 46430.223  # ALC!J83: 52150.75 -> 46430.223
$ws.Cells.Item(83, 11).Value = 117870.75  # ALC!K83: 134281.287 -> 117870.75
$ws.Cells.Item(83, 12).Value = 417872.007  # ALC!L83: 469356.75 -> 417872.007
$ws.Cells.Item(83, 13).Value = -112878.75  # ALC!M83: -129289.287 -> -112878.75
$ws.Cells.Item(83, 14).Value = -427856.007  # ALC!N83: -479340.75 -> -427856.007

$ws.Cells.Item(85, 8).Value = 799.6667  # ALC!H85: 2000 -> 799.6667
$ws.Cells.Item(85, 9).Value = 799.6667  # ALC!I85: 2000 -> 799.6667
$ws.Cells.Item(85, 11).Value = 2399.0001  # ALC!K85: 6000 -> 2399.0001
$ws.Cells.Item(85, 13).Value = -995.0001000000002  # ALC!M85: -4596 -> -995.0001000000002

$ws.Cells.Item(92, 8).Value = 1457.2222  # ALC!H92: 1527.625 -> 1457.2222
$ws.Cells.Item(92, 9).Value = 937.5  # ALC!I92: 946.2 -> 937.5
$ws.Cells.Item(92, 11).Value = 937.5  # ALC!K92: 946.2 -> 937.5
$ws.Cells.Item(92, 13).Value = 310.5  # ALC!M92: 301.8 -> 310.5

$ws.Cells.Item(98, 8).Value = 4740.4546  # ALC!H98: 4687.1816 -> 4740.4546
$ws.Cells.Item(98, 9).Value = 5571.6665  # ALC!I98: 5055.9 -> 5571.6665
$ws.Cells.Item(98, 11).Value = 5571.6665  # ALC!K98: 5055.9 -> 5571.6665
$ws.Cells.Item(98, 13).Value = -4073.6665  # ALC!M98: -3557.9 -> -4073.6665

$ws.Cells.Item(106, 8).Value = 1005  # ALC!H106: 1000000 -> 1005
$ws.Cells.Item(106, 9).Value = 1005  # ALC!I106: 1000000 -> 1005
$ws.Cells.Item(106, 11).Value = 1005  # ALC!K106: 1000000 -> 1005
$ws.Cells.Item(106, 13).Value = -374  # ALC!M106: -999369 -> -374

$ws.Cells.Item(113, 9).Value = 158732020  # ALC!I113: 138890780 -> 158732020
$ws.Cells.Item(113, 10).Value = 41674216  # ALC!J113: 45462590 -> 41674216
$ws.Cells.Item(113, 11).Value = 158732020  # ALC!K113: 138890780 -> 158732020
$ws.Cells.Item(113, 12).Value = 41674216  # ALC!L113: 45462590 -> 41674216
$ws.Cells.Item(113, 13).Value = -158728766  # ALC!M113: -138887526 -> -158728766
$ws.Cells.Item(113, 14).Value = -41680724  # ALC!N113: -45469098 -> -41680724

$ws.Cells.Item(122, 8).Value = 4740.4546  # ALC!H122: 4687.1816 -> 4740.4546
$ws.Cells.Item(122, 9).Value = 5571.6665  # ALC!I122: 5055.9 -> 5571.6665
$ws.Cells.Item(122, 11).Value = 16714.9995  # ALC!K122: 15167.7 -> 16714.9995
$ws.Cells.Item(122, 13).Value = -14264.9995  # ALC!M122: -12717.7 -> -14264.9995

$ws.Cells.Item(132, 8).Value = 1393.7567  # ALC!H132: 1477.7646 -> 1393.7567
$ws.Cells.Item(132, 9).Value = 1410.8235  # ALC!I132: 1482.4375 -> 1410.8235
$ws.Cells.Item(132, 10).Value = 1200.3334  # ALC!J132: 1403 -> 1200.3334
$ws.Cells.Item(132, 11).Value = 4232.470499999999  # ALC!K132: 4447.3125 -> 4232.470499999999
$ws.Cells.Item(132, 12).Value = 3601.0002  # ALC!L132: 4209 -> 3601.0002
$ws.Cells.Item(132, 13).Value = -1702.470499999999  # ALC!M132: -1917.3125 -> -1702.470499999999
$ws.Cells.Item(132, 14).Value = -8661.0002  # ALC!N132: -9269 -> -8661.0002

$ws.Cells.Item(138, 8).Value = 2460.2273  # ALC!H138: 2470.2856 -> 2460.2273
$ws.Cells.Item(138, 10).Value = 2499.9092  # ALC!J138: 2525 -> 2499.9092
$ws.Cells.Item(138, 12).Value = 7499.7276  # ALC!L138: 7575 -> 7499.7276
$ws.Cells.Item(138, 14).Value = -17779.7276  # ALC!N138: -17855 -> -17779.7276

$ws.Cells.Item(141, 8).Value = 2890.125  # ALC!H141: 23812370 -> 2890.125
$ws.Cells.Item(141, 9).Value = 2890.125  # ALC!I141: 23812370 -> 2890.125
$ws.Cells.Item(141, 11).Value = 8670.375  # ALC!K141: 71437110 -> 8670.375
$ws.Cells.Item(141, 13).Value = -3490.375  # ALC!M141: -71431930 -> -3490.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5000.7383  # ARM!H32: 5072.4844 -> 5000.7383
$ws.Cells.Item(32, 9).Value = 4945.18  # ARM!I32: 5020.783 -> 4945.18
$ws.Cells.Item(32, 11).Value = 4945.18  # ARM!K32: 5020.783 -> 4945.18
$ws.Cells.Item(32, 13).Value = -4658.18  # ARM!M32: -4733.783 -> -4658.18

$ws.Cells.Item(61, 8).Value = 3253.26  # ARM!H61: 3303.8367 -> 3253.26
$ws.Cells.Item(61, 9).Value = 2185.6  # ARM!I61: 2185.975 -> 2185.6
$ws.Cells.Item(61, 10).Value = 7523.9  # ARM!J61: 8272.111000000001 -> 7523.9
$ws.Cells.Item(61, 11).Value = 2185.6  # ARM!K61: 2185.975 -> 2185.6
$ws.Cells.Item(61, 12).Value = 7523.9  # ARM!L61: 8272.111000000001 -> 7523.9
$ws.Cells.Item(61, 13).Value = -1973.6  # ARM!M61: -1973.975 -> -1973.6
$ws.Cells.Item(61, 14).Value = -7947.9  # ARM!N61: -8696.111000000001 -> -7947.9

$ws.Cells.Item(74, 8).Value = 40651.605  # ARM!H74: 42625.12 -> 40651.605
$ws.Cells.Item(74, 9).Value = 59600.32  # ARM!I74: 64170 -> 59600.32
$ws.Cells.Item(74, 11).Value = 59600.32  # ARM!K74: 64170 -> 59600.32
$ws.Cells.Item(74, 13).Value = -58726.32  # ARM!M74: -63296 -> -58726.32

$ws.Cells.Item(77, 8).Value = 40651.605  # ARM!H77: 42625.12 -> 40651.605
$ws.Cells.Item(77, 9).Value = 59600.32  # ARM!I77: 64170 -> 59600.32
$ws.Cells.Item(77, 11).Value = 298001.6  # ARM!K77: 320850 -> 298001.6
$ws.Cells.Item(77, 13).Value = -293633.6  # ARM!M77: -316482 -> -293633.6

$ws.Cells.Item(126, 8).Value = 5241.875  # ARM!H126: 5241.7 -> 5241.875
$ws.Cells.Item(126, 9).Value = 5241.875  # ARM!I126: 5241.7 -> 5241.875
$ws.Cells.Item(126, 11).Value = 15725.625  # ARM!K126: 15725.1 -> 15725.625
$ws.Cells.Item(126, 13).Value = -13255.625  # ARM!M126: -13255.1 -> -13255.625

$ws.Cells.Item(136, 8).Value = 3253.26  # ARM!H136: 3303.8367 -> 3253.26
$ws.Cells.Item(136, 9).Value = 2185.6  # ARM!I136: 2185.975 -> 2185.6
$ws.Cells.Item(136, 10).Value = 7523.9  # ARM!J136: 8272.111000000001 -> 7523.9
$ws.Cells.Item(136, 11).Value = 6556.799999999999  # ARM!K136: 6557.924999999999 -> 6556.799999999999
$ws.Cells.Item(136, 12).Value = 22571.7  # ARM!L136: 24816.333 -> 22571.7
$ws.Cells.Item(136, 13).Value = -4006.799999999999  # ARM!M136: -4007.924999999999 -> -4006.799999999999
$ws.Cells.Item(136, 14).Value = -27671.7  # ARM!N136: -29916.333 -> -27671.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3610.4075  # CRP!H16: 3498.6667 -> 3610.4075
$ws.Cells.Item(16, 9).Value = 2206.2942  # CRP!I16: 2142 -> 2206.2942
$ws.Cells.Item(16, 10).Value = 5997.4  # CRP!J16: 5533.6665 -> 5997.4
$ws.Cells.Item(16, 11).Value = 2206.2942  # CRP!K16: 2142 -> 2206.2942
$ws.Cells.Item(16, 12).Value = 5997.4  # CRP!L16: 5533.6665 -> 5997.4
$ws.Cells.Item(16, 13).Value = -1919.2942  # CRP!M16: -1855 -> -1919.2942
$ws.Cells.Item(16, 14).Value = -6571.4  # CRP!N16: -6107.6665 -> -6571.4

$ws.Cells.Item(76, 8).Value = 4953.857  # CRP!H76: 4953.75 -> 4953.857
$ws.Cells.Item(76, 9).Value = 4953.857  # CRP!I76: 4953.75 -> 4953.857
$ws.Cells.Item(76, 11).Value = 4953.857  # CRP!K76: 4953.75 -> 4953.857
$ws.Cells.Item(76, 13).Value = -4638.857  # CRP!M76: -4638.75 -> -4638.857

$ws.Cells.Item(79, 8).Value = 4953.857  # CRP!H79: 4953.75 -> 4953.857
$ws.Cells.Item(79, 9).Value = 4953.857  # CRP!I79: 4953.75 -> 4953.857
$ws.Cells.Item(79, 11).Value = 4953.857  # CRP!K79: 4953.75 -> 4953.857
$ws.Cells.Item(79, 13).Value = -3861.857  # CRP!M79: -3861.75 -> -3861.857

$ws.Cells.Item(113, 8).Value = 3610.4075  # CRP!H113: 3498.6667 -> 3610.4075
$ws.Cells.Item(113, 9).Value = 2206.2942  # CRP!I113: 2142 -> 2206.2942
$ws.Cells.Item(113, 10).Value = 5997.4  # CRP!J113: 5533.6665 -> 5997.4
$ws.Cells.Item(113, 11).Value = 2206.2942  # CRP!K113: 2142 -> 2206.2942
$ws.Cells.Item(113, 12).Value = 5997.4  # CRP!L113: 5533.6665 -> 5997.4
$ws.Cells.Item(113, 13).Value = -36.29419999999982  # CRP!M113: 28 -> -36.29419999999982
$ws.Cells.Item(113, 14).Value = -10337.4  # CRP!N113: -9873.666499999999 -> -10337.4

$ws.Cells.Item(132, 8).Value = 10005283  # CRP!H132: 10261804 -> 10005283
$ws.Cells.Item(132, 9).Value = 2329.0417  # CRP!I132: 2387.4348 -> 2329.0417
$ws.Cells.Item(132, 11).Value = 6987.125100000001  # CRP!K132: 7162.3044 -> 6987.125100000001
$ws.Cells.Item(132, 13).Value = -4457.125100000001  # CRP!M132: -4632.3044 -> -4457.125100000001

$ws.Cells.Item(134, 8).Value = 6358.0625  # CRP!H134: 7100.7144 -> 6358.0625
$ws.Cells.Item(134, 9).Value = 1454.6875  # CRP!I134: 1553.0834 -> 1454.6875
$ws.Cells.Item(134, 11).Value = 4364.0625  # CRP!K134: 4659.2502 -> 4364.0625
$ws.Cells.Item(134, 13).Value = -1829.0625  # CRP!M134: -2124.2502 -> -1829.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 87440  # CUL!H2: 88894.42 -> 87440
$ws.Cells.Item(2, 9).Value = 61.6  # CUL!I2: 15057.45 -> 61.6
$ws.Cells.Item(2, 10).Value = 251274.5  # CUL!J2: 335017.66 -> 251274.5
$ws.Cells.Item(2, 11).Value = 369.6  # CUL!K2: 90344.70000000001 -> 369.6
$ws.Cells.Item(2, 12).Value = 1507647  # CUL!L2: 2010105.96 -> 1507647
$ws.Cells.Item(2, 13).Value = -256.6  # CUL!M2: -90231.70000000001 -> -256.6
$ws.Cells.Item(2, 14).Value = -1507873  # CUL!N2: -2010331.96 -> -1507873

$ws.Cells.Item(11, 8).Value = 1046697.06  # CUL!H11: 1278608.8 -> 1046697.06
$ws.Cells.Item(11, 9).Value = 1278608.6  # CUL!I11: 1438348.6 -> 1278608.6
$ws.Cells.Item(11, 10).Value = 3095  # CUL!J11: 690 -> 3095
$ws.Cells.Item(11, 11).Value = 3835825.8  # CUL!K11: 4315045.800000001 -> 3835825.8
$ws.Cells.Item(11, 12).Value = 9285  # CUL!L11: 2070 -> 9285
$ws.Cells.Item(11, 13).Value = -3835685.8  # CUL!M11: -4314905.800000001 -> -3835685.8
$ws.Cells.Item(11, 14).Value = -9565  # CUL!N11: -2350 -> -9565

$ws.Cells.Item(81, 8).Value = 7286.5  # CUL!H81: 6111.6665 -> 7286.5
$ws.Cells.Item(81, 9).Value = 9006  # CUL!I81: 4503.25 -> 9006
$ws.Cells.Item(81, 10).Value = 6999.9165  # CUL!J81: 6571.2144 -> 6999.9165
$ws.Cells.Item(81, 11).Value = 27018  # CUL!K81: 13509.75 -> 27018
$ws.Cells.Item(81, 12).Value = 20999.7495  # CUL!L81: 19713.6432 -> 20999.7495
$ws.Cells.Item(81, 13).Value = -25895  # CUL!M81: -12386.75 -> -25895
$ws.Cells.Item(81, 14).Value = -23245.7495  # CUL!N81: -21959.6432 -> -23245.7495

$ws.Cells.Item(84, 8).Value = 7286.5  # CUL!H84: 6111.6665 -> 7286.5
$ws.Cells.Item(84, 9).Value = 9006  # CUL!I84: 4503.25 -> 9006
$ws.Cells.Item(84, 10).Value = 6999.9165  # CUL!J84: 6571.2144 -> 6999.9165
$ws.Cells.Item(84, 11).Value = 81054  # CUL!K84: 40529.25 -> 81054
$ws.Cells.Item(84, 12).Value = 62999.2485  # CUL!L84: 59140.9296 -> 62999.2485
$ws.Cells.Item(84, 13).Value = -75438  # CUL!M84: -34913.25 -> -75438
$ws.Cells.Item(84, 14).Value = -74231.2485  # CUL!N84: -70372.9296 -> -74231.2485

$ws.Cells.Item(131, 8).Value = 972.13635  # CUL!H131: 958.087 -> 972.13635
$ws.Cells.Item(131, 9).Value = 747.1053000000001  # CUL!I131: 742.2 -> 747.1053000000001
$ws.Cells.Item(131, 11).Value = 2241.3159  # CUL!K131: 2226.6 -> 2241.3159
$ws.Cells.Item(131, 13).Value = 2798.6841  # CUL!M131: 2813.4 -> 2798.6841

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1851.2122  # GSM!H132: 1908.6774 -> 1851.2122
$ws.Cells.Item(132, 9).Value = 1510.1538  # GSM!I132: 1555.9584 -> 1510.1538
$ws.Cells.Item(132, 11).Value = 4530.4614  # GSM!K132: 4667.8752 -> 4530.4614
$ws.Cells.Item(132, 13).Value = -2000.4614  # GSM!M132: -2137.8752 -> -2000.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1098.2858  # LTW!H22: 1126.1765 -> 1098.2858
$ws.Cells.Item(22, 10).Value = 2562  # LTW!J22: 2906.5715 -> 2562
$ws.Cells.Item(22, 12).Value = 2562  # LTW!L22: 2906.5715 -> 2562
$ws.Cells.Item(22, 14).Value = -3152  # LTW!N22: -3496.5715 -> -3152

$ws.Cells.Item(27, 8).Value = 1098.2858  # LTW!H27: 1126.1765 -> 1098.2858
$ws.Cells.Item(27, 10).Value = 2562  # LTW!J27: 2906.5715 -> 2562
$ws.Cells.Item(27, 12).Value = 2562  # LTW!L27: 2906.5715 -> 2562
$ws.Cells.Item(27, 14).Value = -2776  # LTW!N27: -3120.5715 -> -2776

$ws.Cells.Item(33, 8).Value = 16674332  # LTW!H33: 12506971 -> 16674332
$ws.Cells.Item(33, 9).Value = 16674332  # LTW!I33: 12506971 -> 16674332
$ws.Cells.Item(33, 11).Value = 16674332  # LTW!K33: 12506971 -> 16674332
$ws.Cells.Item(33, 13).Value = -16674042  # LTW!M33: -12506681 -> -16674042

$ws.Cells.Item(34, 8).Value = 1684.3334  # LTW!H34: 1696.75 -> 1684.3334
$ws.Cells.Item(34, 9).Value = 1684.3334  # LTW!I34: 1696.75 -> 1684.3334
$ws.Cells.Item(34, 11).Value = 1684.3334  # LTW!K34: 1696.75 -> 1684.3334
$ws.Cells.Item(34, 13).Value = -1512.3334  # LTW!M34: -1524.75 -> -1512.3334

$ws.Cells.Item(39, 8).Value = 0  # LTW!H39: 2000 -> 0
$ws.Cells.Item(39, 9).Value = 0  # LTW!I39: 2000 -> 0
$ws.Cells.Item(39, 11).Value = 0  # LTW!K39: 2000 -> 0
$ws.Cells.Item(39, 13).ClearContents()  # LTW!M39: -1540 -> (removed)

$ws.Cells.Item(93, 8).Value = 5291.625  # LTW!H93: 5618.067 -> 5291.625
$ws.Cells.Item(93, 10).Value = 8616.5  # LTW!J93: 10260.8 -> 8616.5
$ws.Cells.Item(93, 12).Value = 8616.5  # LTW!L93: 10260.8 -> 8616.5
$ws.Cells.Item(93, 14).Value = -11112.5  # LTW!N93: -12756.8 -> -11112.5

$ws.Cells.Item(116, 8).Value = 56958  # LTW!H116: 0 -> 56958
$ws.Cells.Item(116, 10).Value = 56958  # LTW!J116: 0 -> 56958
$ws.Cells.Item(116, 12).Value = 56958  # LTW!L116: 0 -> 56958
$ws.Cells.Item(116, 14).Value = -66136  # LTW!N116: None -> -66136

$ws.Cells.Item(132, 8).Value = 4768.75  # LTW!H132: 4933 -> 4768.75
$ws.Cells.Item(132, 9).Value = 2830.0334  # LTW!I132: 2961.3928 -> 2830.0334
$ws.Cells.Item(132, 11).Value = 8490.100199999999  # LTW!K132: 8884.178400000001 -> 8490.100199999999
$ws.Cells.Item(132, 13).Value = -5960.100199999999  # LTW!M132: -6354.178400000001 -> -5960.100199999999

$ws.Cells.Item(136, 8).Value = 7516.0557  # LTW!H136: 7895.125 -> 7516.0557
$ws.Cells.Item(136, 9).Value = 3572.25  # LTW!I136: 2661 -> 3572.25
$ws.Cells.Item(136, 11).Value = 10716.75  # LTW!K136: 7983 -> 10716.75
$ws.Cells.Item(136, 13).Value = -8166.75  # LTW!M136: -5433 -> -8166.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1322.3334  # WVR!H100: 1259.0769 -> 1322.3334
$ws.Cells.Item(100, 9).Value = 1205.75  # WVR!I100: 1064.6 -> 1205.75
$ws.Cells.Item(100, 11).Value = 2411.5  # WVR!K100: 2129.2 -> 2411.5
$ws.Cells.Item(100, 13).Value = -1870.5  # WVR!M100: -1588.2 -> -1870.5

$ws.Cells.Item(122, 8).Value = 152477.81  # WVR!H122: 105721.8 -> 152477.81
$ws.Cells.Item(122, 9).Value = 224314.78  # WVR!I122: 139790.2 -> 224314.78
$ws.Cells.Item(122, 10).Value = 8803.888999999999  # WVR!J122: 6923.4 -> 8803.888999999999
$ws.Cells.Item(122, 11).Value = 672944.34  # WVR!K122: 419370.6 -> 672944.34
$ws.Cells.Item(122, 12).Value = 26411.667  # WVR!L122: 20770.2 -> 26411.667
$ws.Cells.Item(122, 13).Value = -670494.34  # WVR!M122: -416920.6 -> -670494.34
$ws.Cells.Item(122, 14).Value = -31311.667  # WVR!N122: -25670.2 -> -31311.667

$ws.Cells.Item(124, 8).Value = 34506  # WVR!H124: 33746.25 -> 34506
$ws.Cells.Item(124, 9).Value = 0  # WVR!I124: 23390 -> 0
$ws.Cells.Item(124, 10).Value = 34506  # WVR!J124: 37198.332 -> 34506
$ws.Cells.Item(124, 11).Value = 0  # WVR!K124: 23390 -> 0
$ws.Cells.Item(124, 12).Value = 34506  # WVR!L124: 37198.332 -> 34506
$ws.Cells.Item(124, 13).ClearContents()  # WVR!M124: -18480 -> (removed)
$ws.Cells.Item(124, 14).Value = -44326  # WVR!N124: -47018.332 -> -44326

$ws.Cells.Item(132, 8).Value = 15350.63  # WVR!H132: 13833.767 -> 15350.63
$ws.Cells.Item(132, 9).Value = 8303.866  # WVR!I132: 7335.5293 -> 8303.866
$ws.Cells.Item(132, 10).Value = 24159.084  # WVR!J132: 22331.46 -> 24159.084
$ws.Cells.Item(132, 11).Value = 24911.598  # WVR!K132: 22006.5879 -> 24911.598
$ws.Cells.Item(132, 12).Value = 72477.25199999999  # WVR!L132: 66994.38 -> 72477.25199999999
$ws.Cells.Item(132, 13).Value = -22381.598  # WVR!M132: -19476.5879 -> -22381.598
$ws.Cells.Item(132, 14).Value = -77537.25199999999  # WVR!N132: -72054.38 -> -77537.25199999999
